$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(population census results)" subtitle text from A2 - it becomes a
# plain empty row once the content + explicit formatting are both cleared.
$ws.Range("A2").ClearContents()
$ws.Range("A2").Style = "Normal"

# The row that used to hold only a formatting placeholder under A2 (old row 3)
# is removed entirely, shifting everything below it up by one row.
$ws.Rows(3).Delete()

# The 1989 and 2002 columns are dropped, leaving only the 2014 column (which
# becomes column B).
$ws.Range("B:C").Delete()

# Re-select A2, matching the saved selection in the final sheet.
$ws.Range("A2").Select()

# Rename the sheet from the generic "1" to the municipality's name.
$ws.Name = "ლაგოდეხი"
